# Generate Report for Handoff
#
# - Status moves from "Handed back: in sync with en-US" to "Ready for
#   handoff", with refreshed handoff timestamps.
# - The e464a940-... entry is dropped from this handoff run, so its row
#   (row 3) is removed from the Overview / zh-cn / de-de sheets, together
#   with the hyperlinks that pointed at it.
#
# NB: this engine's `Range(...).Hyperlinks.Delete()` clears every
# hyperlink on the worksheet (not just the targeted cell), so the
# reliable way to drop only the stale ones is: wipe the sheet's
# hyperlinks once, then re-add the ones that should survive with their
# original target URL + display text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-33-11 18:33:45"

$ov.Range("A3:D3").EntireRow.Delete()

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md") | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-11 18:33:42"

$zh.Range("A3:K3").EntireRow.Delete()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e08bd87dba2db99d32b6aeee3d484b80778344a2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ad440bb342af6d37a87246565b00caa34c353763/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d609d04d18b7a49195f24d4e556069989c8621c5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-11 18:33:45"

$de.Range("A3:K3").EntireRow.Delete()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/487f785d7c0172c05129dc0b4f790bed39accb9e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/048aa54ae6b221acb9ae568f084b34ebda184e02/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dfa18673c9cd7f0d723c813087ae62954778583b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf") | Out-Null
